# Update LR-pair TPM values (Areg-Egfr) per new TPM computation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2247043333333333
$ws.Range("H2").Value = 0.674113
$ws.Range("I2").Value = 0.2389319335355998
$ws.Range("J2").Value = 0.2389319335355999
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 0.3824309711285555
$ws.Range("R2").Value = 3.441878740157
$ws.Range("S2").Value = 0.005031005477199709
$ws.Range("T2").Value = 0.00503100547719971
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2247043333333333
$ws.Range("H3").Value = 0.674113
$ws.Range("I3").Value = 0.2389319335355998
$ws.Range("J3").Value = 0.2389319335355999
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 14.044907741337
$ws.Range("R3").Value = 126.404169672033
$ws.Range("S3").Value = 0.1847653906400756
$ws.Range("T3").Value = 0.1847653906400757
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2247043333333333
$ws.Range("H4").Value = 0.674113
$ws.Range("I4").Value = 0.2389319335355998
$ws.Range("J4").Value = 0.2389319335355999
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 0.097492244399
$ws.Range("R4").Value = 0.8774301995910001
$ws.Range("S4").Value = 0.001282542609215046
$ws.Range("T4").Value = 0.001282542609215047
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2247043333333333
$ws.Range("H5").Value = 0.674113
$ws.Range("I5").Value = 0.2389319335355998
$ws.Range("J5").Value = 0.2389319335355999
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 3.589063898464
$ws.Range("R5").Value = 32.301575086176
$ws.Range("S5").Value = 0.04721531856561464
$ws.Range("T5").Value = 0.04721531856561465
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2247043333333333
$ws.Range("H6").Value = 0.674113
$ws.Range("I6").Value = 0.2389319335355998
$ws.Range("J6").Value = 0.2389319335355999
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 0.04847284427944444
$ws.Range("R6").Value = 0.436255598515
$ws.Range("S6").Value = 0.0006376762434948209
$ws.Range("T6").Value = 0.0006376762434948211
$ws.Range("G7").Value = 0.4451493333333333
$ws.Range("H7").Value = 1.335448
$ws.Range("I7").Value = 0.4733348456063742
$ws.Range("J7").Value = 0.4733348456063743
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 0.7576128564968888
$ws.Range("R7").Value = 6.818515708472
$ws.Range("S7").Value = 0.009966646841872798
$ws.Range("T7").Value = 0.0099666468418728
$ws.Range("G8").Value = 0.4451493333333333
$ws.Range("H8").Value = 1.335448
$ws.Range("I8").Value = 0.4733348456063742
$ws.Range("J8").Value = 0.4733348456063743
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("Q8").Value = 27.823590337752
$ws.Range("R8").Value = 250.412313039768
$ws.Range("S8").Value = 0.3660285017489764
$ws.Range("T8").Value = 0.3660285017489765
$ws.Range("G9").Value = 0.4451493333333333
$ws.Range("H9").Value = 1.335448
$ws.Range("I9").Value = 0.4733348456063742
$ws.Range("J9").Value = 0.4733348456063743
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 0.193136496104
$ws.Range("R9").Value = 1.738228464936
$ws.Range("S9").Value = 0.002540774265428816
$ws.Range("T9").Value = 0.002540774265428817
$ws.Range("G10").Value = 0.4451493333333333
$ws.Range("H10").Value = 1.335448
$ws.Range("I10").Value = 0.4733348456063742
$ws.Range("J10").Value = 0.4733348456063743
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 7.110096089344
$ws.Range("R10").Value = 63.990864804096
$ws.Range("S10").Value = 0.09353565759422076
$ws.Range("T10").Value = 0.09353565759422079
$ws.Range("G11").Value = 0.4451493333333333
$ws.Range("H11").Value = 1.335448
$ws.Range("I11").Value = 0.4733348456063742
$ws.Range("J11").Value = 0.4733348456063743
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 0.09602687227111112
$ws.Range("R11").Value = 0.86424185044
$ws.Range("S11").Value = 0.001263265155875457
$ws.Range("T11").Value = 0.001263265155875457
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2705996666666666
$ws.Range("H12").Value = 0.8117989999999999
$ws.Range("I12").Value = 0.2877332208580259
$ws.Range("J12").Value = 0.2877332208580259
$ws.Range("M12").Value = 1.701929666666667
$ws.Range("N12").Value = 5.105789
$ws.Range("O12").Value = 0.02105622887134972
$ws.Range("P12").Value = 0.02105622887134972
$ws.Range("Q12").Value = 0.460541600490111
$ws.Range("R12").Value = 4.144874404411
$ws.Range("S12").Value = 0.00605857655227721
$ws.Range("T12").Value = 0.00605857655227721
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2705996666666666
$ws.Range("H13").Value = 0.8117989999999999
$ws.Range("I13").Value = 0.2877332208580259
$ws.Range("J13").Value = 0.2877332208580259
$ws.Range("O13").Value = 0.7732971809418951
$ws.Range("P13").Value = 0.7732971809418953
$ws.Range("Q13").Value = 16.913547223551
$ws.Range("R13").Value = 152.221925011959
$ws.Range("S13").Value = 0.2225032885528431
$ws.Range("T13").Value = 0.2225032885528432
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2705996666666666
$ws.Range("H14").Value = 0.8117989999999999
$ws.Range("I14").Value = 0.2877332208580259
$ws.Range("J14").Value = 0.2877332208580259
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4338690000000001
$ws.Range("N14").Value = 1.301607
$ws.Range("O14").Value = 0.005367815805265532
$ws.Range("P14").Value = 0.005367815805265533
$ws.Range("Q14").Value = 0.117404806777
$ws.Range("R14").Value = 1.056643260993
$ws.Range("S14").Value = 0.001544498930621669
$ws.Range("T14").Value = 0.00154449893062167
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2705996666666666
$ws.Range("H15").Value = 0.8117989999999999
$ws.Range("I15").Value = 0.2877332208580259
$ws.Range("J15").Value = 0.2877332208580259
$ws.Range("M15").Value = 15.972384
$ws.Range("N15").Value = 47.917152
$ws.Range("O15").Value = 0.1976099128607259
$ws.Range("P15").Value = 0.1976099128607259
$ws.Range("Q15").Value = 4.322121786272
$ws.Range("R15").Value = 38.899096076448
$ws.Range("S15").Value = 0.0568589367008905
$ws.Range("T15").Value = 0.05685893670089052
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2705996666666666
$ws.Range("H16").Value = 0.8117989999999999
$ws.Range("I16").Value = 0.2877332208580259
$ws.Range("J16").Value = 0.2877332208580259
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2157183333333333
$ws.Range("N16").Value = 0.647155
$ws.Range("O16").Value = 0.002668861520763652
$ws.Range("P16").Value = 0.002668861520763652
$ws.Range("Q16").Value = 0.05837330909388888
$ws.Range("R16").Value = 0.525359781845
$ws.Range("S16").Value = 0.0007679201213933749
$ws.Range("T16").Value = 0.000767920121393375
